# Q2 is proper good
# The sample-name header cells (row 1, columns B:DQ) and the mirrored
# sample-name cells in column A (rows 2:121) were stored as the string
# representation of a 1-tuple, e.g. "('ANG_1',)". This cleans them up to
# plain strings, e.g. "ANG_1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Strip-TupleRepr($val) {
    if ($null -eq $val) {
        return $val
    }
    if (($val.StartsWith("('")) -and ($val.EndsWith("',)"))) {
        return $val.Substring(2, $val.Length - 5)
    }
    return $val
}

# Row 1 header labels, columns B (2) through DQ (121)
for ($col = 2; $col -le 121; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    $cell.Value = Strip-TupleRepr $val
}

# Column A mirrored labels, rows 2 through 121
for ($row = 2; $row -le 121; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $val = $cell.Value()
    $cell.Value = Strip-TupleRepr $val
}
